# Word COM-interop script implementing the "week 9" edit described by the diff.
#
# Substance of the edit (once the cosmetic proofErr/run-splitting noise that
# Word's spell/grammar checker adds is stripped away) is a series of
# "we"/"We" -> "I" wording fixes throughout the write-up, plus relocating the
# auto-managed "_GoBack" bookmark from its old position to right after the
# "Overview" heading run (i.e. the last place the author actually edited).

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARNING: find/replace failed for: $old"
    }
}

# --- Wording fixes ("we"/"We" -> "I") -----------------------------------

Replace-Text "After creating the account, we get this AgentKey" "After creating the account, I get this AgentKey"

Replace-Text "when we add this key to the program" "when I add this key to the program"

Replace-Text "One conclusion we can draw:" "One conclusion I can draw:"

Replace-Text "it is a great tool to find out the result when we use the pprof." "it is a great tool to find out the result when I use the pprof."

Replace-Text "We can see that after doing some modification" "I can see that after doing some modification"

Replace-Text "However, we can combine this two tools together" "However, I can combine this two tools together"

Replace-Text "We can see that in memory profiling dashboard, it shows that we use the func3" "I can see that in memory profiling dashboard, it shows that I use the func3"

Replace-Text "after running the na$([char]0xEF)ve program, we can see the top10 list and after checking the svg file, we found out that" "after running the na$([char]0xEF)ve program, I can see the top10 list and after checking the svg file, I found out that"

Replace-Text "We go deeper and check what getStatesTags do, and find out that it calls os.Hostname(), and we then realized" "I go deeper and check what getStatesTags do, and find out that it calls os.Hostname(), and I then realized"

Replace-Text ". We then cache this hostname since it$([char]0x2019)s same for one request. We store the hostname as a global variable" ". I then cache this hostname since it$([char]0x2019)s same for one request. I store the hostname as a global variable"

Replace-Text "Based on all the above experiment and research, we should use stackImpect" "Based on all the above experiment and research, I should use stackImpect"

Replace-Text "We used couple days to research and realized that stackImpect is powerful(maybe), but is not a helpful tool for our project. We can use it as a $([char]0x201C)monitor$([char]0x201D) though." "I used couple days to research and realized that stackImpect is powerful(maybe), but is not a helpful tool for our project. I can use it as a $([char]0x201C)monitor$([char]0x201D) though."

# --- Relocate the "_GoBack" bookmark ------------------------------------
# It currently sits right before the closing ")" after "...drop down
# immediately". It should end up as a zero-length bookmark right after the
# "Overview" run. Directly collapsing a Range to the exact end-of-paragraph
# boundary and calling Bookmarks.Add there mis-places the markers in this
# runtime, so we sidestep that boundary case by inserting a one-character
# placeholder at the target spot, bookmarking around it, then deleting the
# placeholder (the bookmark collapses cleanly to the right spot).

$existing = $d.Bookmarks("_GoBack")
if ($existing -ne $null) {
    $existing.Delete()
}

$findRng = $d.Content
$found = $findRng.Find.Execute("Overview")
if ($found) {
    $insRng = $findRng.Duplicate
    $insRng.Collapse(0)
    $insStart = $insRng.Start
    $insRng.InsertAfter("X")
    $d.Bookmarks.Add("_GoBack", $insRng)
    $delRng = $d.Range($insStart, $insStart + 1)
    $delRng.Text = ""
} else {
    Write-Output "WARNING: could not find Overview heading to place _GoBack bookmark"
}

Write-Output "done"
